{"js": "// Fill in the \"Name Sign-up\" column of the driver-clinic schedule table\n// with the drivers who completed the session (diff: adding completed\n// driver clinic 1 names to the 2023 TFR agenda).\n//\n// Table layout (0-indexed rows):\n//   row 6  -> 11:00am \u2013 11:30am | Sim Driver 1 | (Name Sign-up)\n//   row 7  -> 11:30am \u2013 12:00pm | Sim Driver 2 | (Name Sign-up)\n//   row 8  -> 12:00pm \u2013 12:30pm | Sim Driver 3 | (Name Sign-up)\n//   row 9  -> 12:30pm \u2013 1:00pm  | Sim Driver 4 | (Name Sign-up)\n//   row 13 -> 1:30pm \u2013 2:00pm   | Sim Driver 1 | (Name Sign-up)\n//   row 14 -> 2:00pm \u2013 2:30pm   | Sim Driver 2 | (Name Sign-up)\n//   row 15 -> 2:30pm \u2013 3:00pm   | Sim Driver 3 | (Name Sign-up)\n//   row 16 -> 3:00pm \u2013 3:30pm   | Sim Driver 4 | (Name Sign-up)\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, nameSignUpColumnIndex, name]\nconst updates = [\n  [6, 2, \"Bern\"],\n  [7, 2, \"Anne Marie\"],\n  [8, 2, \"Kaci\"],\n  [9, 2, \"Anne Marie\"],\n  [13, 2, \"Fred\"],\n  [14, 2, \"Bern\"],\n  [15, 2, \"Fred\"],\n  [16, 2, \"Kaci\"],\n];\n\nfor (const [rowIndex, colIndex, name] of updates) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  // The target cell's lone (empty) paragraph gets the driver's name,\n  // preserving the paragraph's existing formatting (centered).\n  paragraphs.items[0].insertText(name, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fill in the \"Name Sign-up\" column of the driver-clinic schedule table\n# with the drivers who completed the session (diff: adding completed\n# driver clinic 1 names to the 2023 TFR agenda).\n#\n# Table layout (1-indexed rows/cols, as COM Table.Cell expects):\n#   row 7  -> 11:00am - 11:30am | Sim Driver 1 | (Name Sign-up)\n#   row 8  -> 11:30am - 12:00pm | Sim Driver 2 | (Name Sign-up)\n#   row 9  -> 12:00pm - 12:30pm | Sim Driver 3 | (Name Sign-up)\n#   row 10 -> 12:30pm - 1:00pm  | Sim Driver 4 | (Name Sign-up)\n#   row 14 -> 1:30pm - 2:00pm   | Sim Driver 1 | (Name Sign-up)\n#   row 15 -> 2:00pm - 2:30pm   | Sim Driver 2 | (Name Sign-up)\n#   row 16 -> 2:30pm - 3:00pm   | Sim Driver 3 | (Name Sign-up)\n#   row 17 -> 3:00pm - 3:30pm   | Sim Driver 4 | (Name Sign-up)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 7;  Col = 3; Name = \"Bern\" },\n    @{ Row = 8;  Col = 3; Name = \"Anne Marie\" },\n    @{ Row = 9;  Col = 3; Name = \"Kaci\" },\n    @{ Row = 10; Col = 3; Name = \"Anne Marie\" },\n    @{ Row = 14; Col = 3; Name = \"Fred\" },\n    @{ Row = 15; Col = 3; Name = \"Bern\" },\n    @{ Row = 16; Col = 3; Name = \"Fred\" },\n    @{ Row = 17; Col = 3; Name = \"Kaci\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $t.Cell($u.Row, $u.Col)\n    $cell.Range.Text = $u.Name\n}\n"}
